$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2341.1428
$ws.Range("I40").Value = 2493
$ws.Range("J40").Value = 2315.8333
$ws.Range("K40").Value = 2493
$ws.Range("L40").Value = 2315.8333
$ws.Range("M40").Value = -2318
$ws.Range("N40").Value = -2665.8333
$ws.Range("H64").Value = 6194.5835
$ws.Range("J64").Value = 6508.25
$ws.Range("L64").Value = 6508.25
$ws.Range("N64").Value = -7004.25
$ws.Range("H67").Value = 6194.5835
$ws.Range("J67").Value = 6508.25
$ws.Range("L67").Value = 6508.25
$ws.Range("N67").Value = -8224.25
$ws.Range("H76").Value = 18890.1
$ws.Range("I76").Value = 14129
$ws.Range("J76").Value = 29999.334
$ws.Range("K76").Value = 14129
$ws.Range("L76").Value = 29999.334
$ws.Range("M76").Value = -13814
$ws.Range("N76").Value = -30629.334
$ws.Range("H79").Value = 18890.1
$ws.Range("I79").Value = 14129
$ws.Range("J79").Value = 29999.334
$ws.Range("K79").Value = 14129
$ws.Range("L79").Value = 29999.334
$ws.Range("M79").Value = -13037
$ws.Range("N79").Value = -32183.334
$ws.Range("H112").Value = 4688.5083
$ws.Range("I112").Value = 3491
$ws.Range("J112").Value = 4709.1553
$ws.Range("K112").Value = 10473
$ws.Range("L112").Value = 14127.4659
$ws.Range("M112").Value = -9365
$ws.Range("N112").Value = -16343.4659
$ws.Range("H138").Value = 4946.7437
$ws.Range("I138").Value = 3282.3635
$ws.Range("K138").Value = 9847.0905
$ws.Range("M138").Value = -4707.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1216.3914
$ws.Range("I2").Value = 1156.1628
$ws.Range("K2").Value = 1156.1628
$ws.Range("M2").Value = -1043.1628
$ws.Range("H63").Value = 2619.6
$ws.Range("I63").Value = 2867.6667
$ws.Range("J63").Value = 2247.5
$ws.Range("K63").Value = 2867.6667
$ws.Range("L63").Value = 2247.5
$ws.Range("M63").Value = -2181.6667
$ws.Range("N63").Value = -3619.5
$ws.Range("H66").Value = 2619.6
$ws.Range("I66").Value = 2867.6667
$ws.Range("J66").Value = 2247.5
$ws.Range("K66").Value = 14338.3335
$ws.Range("L66").Value = 11237.5
$ws.Range("M66").Value = -10906.3335
$ws.Range("N66").Value = -18101.5
$ws.Range("H116").Value = 1216.3914
$ws.Range("I116").Value = 1156.1628
$ws.Range("K116").Value = 1156.1628
$ws.Range("M116").Value = 1137.8372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1216.3914
$ws.Range("I3").Value = 1156.1628
$ws.Range("K3").Value = 1156.1628
$ws.Range("M3").Value = -1042.1628
$ws.Range("H20").Value = 5733.375
$ws.Range("I20").Value = 5671
$ws.Range("K20").Value = 5671
$ws.Range("M20").Value = -5424
$ws.Range("H86").Value = 2261.0833
$ws.Range("I86").Value = 2169.8572
$ws.Range("J86").Value = 2388.8
$ws.Range("K86").Value = 2169.8572
$ws.Range("L86").Value = 2388.8
$ws.Range("M86").Value = -1046.8572
$ws.Range("N86").Value = -4634.8
$ws.Range("H89").Value = 2261.0833
$ws.Range("I89").Value = 2169.8572
$ws.Range("J89").Value = 2388.8
$ws.Range("K89").Value = 10849.286
$ws.Range("L89").Value = 11944
$ws.Range("M89").Value = -5233.286
$ws.Range("N89").Value = -23176
$ws.Range("H134").Value = 2503.5334
$ws.Range("I134").Value = 2405.6
$ws.Range("K134").Value = 7216.799999999999
$ws.Range("M134").Value = -4681.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5750
$ws.Range("J62").Value = 5750
$ws.Range("L62").Value = 5750
$ws.Range("N62").Value = -6998
$ws.Range("H65").Value = 5750
$ws.Range("J65").Value = 5750
$ws.Range("L65").Value = 28750
$ws.Range("N65").Value = -34990

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 102.75
$ws.Range("J33").Value = 156
$ws.Range("L33").Value = 936
$ws.Range("N33").Value = -1502
$ws.Range("H55").Value = 11285.143
$ws.Range("I55").Value = 4000
$ws.Range("J55").Value = 12499.333
$ws.Range("K55").Value = 12000
$ws.Range("L55").Value = 37497.999
$ws.Range("M55").Value = -11823
$ws.Range("N55").Value = -37851.999
$ws.Range("H132").Value = 664.5
$ws.Range("I132").Value = 787.3333
$ws.Range("J132").Value = 443.4
$ws.Range("K132").Value = 7085.9997
$ws.Range("L132").Value = 3990.6
$ws.Range("M132").Value = -4555.9997
$ws.Range("N132").Value = -9050.6
$ws.Range("H136").Value = 5943.758
$ws.Range("I136").Value = 5271.32
$ws.Range("K136").Value = 15813.96
$ws.Range("M136").Value = -10713.96
$ws.Range("H139").Value = 4027.25
$ws.Range("I139").Value = 3855.5881
$ws.Range("K139").Value = 11566.7643
$ws.Range("M139").Value = -6426.764299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4217.6665
$ws.Range("I70").Value = 4301.4
$ws.Range("J70").Value = 3799
$ws.Range("K70").Value = 4301.4
$ws.Range("L70").Value = 3799
$ws.Range("M70").Value = -4031.4
$ws.Range("N70").Value = -4339
$ws.Range("H73").Value = 4217.6665
$ws.Range("I73").Value = 4301.4
$ws.Range("J73").Value = 3799
$ws.Range("K73").Value = 4301.4
$ws.Range("L73").Value = 3799
$ws.Range("M73").Value = -3365.4
$ws.Range("N73").Value = -5671
$ws.Range("H80").Value = 6280.6523
$ws.Range("J80").Value = 7088.1
$ws.Range("L80").Value = 7088.1
$ws.Range("N80").Value = -9084.1
$ws.Range("H83").Value = 6280.6523
$ws.Range("J83").Value = 7088.1
$ws.Range("L83").Value = 35440.5
$ws.Range("N83").Value = -45424.5
$ws.Range("H102").Value = 934.4286
$ws.Range("I102").Value = 932
$ws.Range("K102").Value = 932
$ws.Range("M102").Value = 690
$ws.Range("H132").Value = 3238.9333
$ws.Range("I132").Value = 3406.0356
$ws.Range("K132").Value = 10218.1068
$ws.Range("M132").Value = -7688.106800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 295551.5
$ws.Range("J5").Value = 295551.5
$ws.Range("L5").Value = 295551.5
$ws.Range("N5").Value = -295775.5
$ws.Range("H126").Value = 4173.778
$ws.Range("I126").Value = 3445.5
$ws.Range("K126").Value = 10336.5
$ws.Range("M126").Value = -7866.5
$ws.Range("H136").Value = 22204.25
$ws.Range("I136").Value = 26669.545
$ws.Range("K136").Value = 80008.63499999999
$ws.Range("M136").Value = -77458.63499999999
